$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.691.02"
$ws.Range("E2").Value = "  -1.47%  "

$ws.Range("D3").Value = "2.906.40"
$ws.Range("E3").Value = "  -2.05%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.26%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.555"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.12%  "

$ws.Range("D9").Value = "2.914.23"
$ws.Range("E9").Value = "  -2.13%  "

$ws.Range("E10").Value = "  -4.69%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.05"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.360"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.78%  "

$ws.Range("D13").Value = "3.416.86"
$ws.Range("E13").Value = "  -2.00%  "

$ws.Range("E14").Value = "  +2.20%  "

$ws.Range("D15").Value = "60.665.41"
$ws.Range("E15").Value = "  -1.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.93%  "

$ws.Range("D17").Value = "2.914.51"
$ws.Range("E17").Value = "  -2.01%  "

$ws.Range("E18").Value = "  -3.79%  "

$ws.Range("E19").Value = "  -2.42%  "

$ws.Range("E20").Value = "  -2.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "361.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.60%  "

$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("E24").Value = "  +0.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.31%  "

$ws.Range("E26").Value = "  -3.03%  "

$ws.Range("E27").Value = "  -3.66%  "

$ws.Range("E28").Value = "  +0.16%  "

$ws.Range("E29").Value = "  -4.97%  "

$ws.Range("E30").Value = "  -9.31%  "

$ws.Range("E32").Value = "  -2.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "151.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.94%  "

$ws.Range("E35").Value = "  -5.83%  "

$ws.Range("E36").Value = "  -6.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.57%  "

$ws.Range("E38").Value = "  -4.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.61%  "

$ws.Range("E40").Value = "  -4.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.90%  "

$ws.Range("D42").Value = "2.295.29"
$ws.Range("E42").Value = "  -4.77%  "

$ws.Range("E43").Value = "  -2.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0587"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0239"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0924"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.17%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "249.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.27%  "
